# Apply Crypto price/volume updates (GitHub Actions refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds prices as text (e.g. '61.701.49' uses '.' as a thousands
# separator, which Excel would otherwise mis-parse as a number). Force the
# whole price column to text first so every write lands as a string, then
# drop back to the default ('Normal') style so we don't leave cells visibly
# reformatted -- matches the source workbook where these cells carry no
# explicit style.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "61.701.49"
$ws.Range("E2").Value = "  +0.42%  "
$ws.Range("D3").Value = "3.424.84"
$ws.Range("E3").Value = "  +1.51%  "
$ws.Range("E4").Value = "  -0.73%  "
$ws.Range("D5").Value = "407.75"
$ws.Range("E5").Value = "  -1.69%  "
$ws.Range("D6").Value = "128.83"
$ws.Range("E6").Value = "  +15.06%  "
$ws.Range("D7").Value = "3.414.25"
$ws.Range("E7").Value = "  +4.19%  "
$ws.Range("D8").Value = "0.601"
$ws.Range("E8").Value = "  +4.51%  "
$ws.Range("E9").Value = "  -0.49%  "
$ws.Range("D10").Value = "0.675"
$ws.Range("E10").Value = "  +6.95%  "
$ws.Range("D11").Value = "0.127"
$ws.Range("E11").Value = "  +16.27%  "
$ws.Range("D12").Value = "42.28"
$ws.Range("E12").Value = "  +7.78%  "
$ws.Range("D13").Value = "0.141"
$ws.Range("E13").Value = "  -0.70%  "
$ws.Range("D14").Value = "3.970.55"
$ws.Range("E14").Value = "  -1.13%  "
$ws.Range("D15").Value = "8.58"
$ws.Range("E15").Value = "  +5.28%  "
$ws.Range("D16").Value = "19.78"
$ws.Range("E16").Value = "  +2.54%  "
$ws.Range("D17").Value = "3.426.06"
$ws.Range("E17").Value = "  -2.32%  "
$ws.Range("B18").Value = "Uniswap"
$ws.Range("C18").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D18").Value = "11.59"
$ws.Range("E18").Value = "  +10.00%  "
$ws.Range("B19").Value = "WrappedBTC"
$ws.Range("C19").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D19").Value = "61.899.76"
$ws.Range("E19").Value = "  -0.86%  "
$ws.Range("E20").Value = "  +0.43%  "
$ws.Range("D21").Value = "0.0000135"
$ws.Range("E21").Value = "  +20.13%  "
$ws.Range("D22").Value = "3.24"
$ws.Range("E22").Value = "  -0.80%  "
$ws.Range("D23").Value = "83.38"
$ws.Range("E23").Value = "  +11.78%  "
$ws.Range("D24").Value = "12.99"
$ws.Range("E24").Value = "  +7.09%  "
$ws.Range("D25").Value = "308.22"
$ws.Range("E25").Value = "  +4.93%  "
$ws.Range("E26").Value = "  -2.40%  "
$ws.Range("D27").Value = "8.57"
$ws.Range("E27").Value = "  +14.85%  "
$ws.Range("D28").Value = "29.73"
$ws.Range("E28").Value = "  +2.85%  "
$ws.Range("E29").Value = "  +0.57%  "
$ws.Range("D30").Value = "7.48"
$ws.Range("E30").Value = "  +2.35%  "
$ws.Range("E31").Value = "  +10.38%  "
$ws.Range("D32").Value = "0.115"
$ws.Range("E32").Value = "  +4.46%  "
$ws.Range("D33").Value = "11.66"
$ws.Range("E33").Value = "  +4.47%  "
$ws.Range("D34").Value = "42.86"
$ws.Range("E34").Value = "  +12.29%  "
$ws.Range("D35").Value = "2.55"
$ws.Range("E35").Value = "  +10.29%  "
$ws.Range("D36").Value = "0.999"
$ws.Range("E36").Value = "  +0.16%  "
$ws.Range("D37").Value = "0.0484"
$ws.Range("E37").Value = "  -1.30%  "
$ws.Range("D38").Value = "52.22"
$ws.Range("E38").Value = "  +0.47%  "
$ws.Range("D39").Value = "1.00"
$ws.Range("E39").Value = "  -0.51%  "
$ws.Range("D40").Value = "3.40"
$ws.Range("E40").Value = "  +3.17%  "
$ws.Range("D41").Value = "3.00"
$ws.Range("E41").Value = "  -1.68%  "
$ws.Range("E42").Value = "  +4.07%  "
$ws.Range("D43").Value = "1.96"
$ws.Range("E43").Value = "  +3.86%  "
$ws.Range("D44").Value = "135.69"
$ws.Range("E44").Value = "  -1.54%  "
$ws.Range("D45").Value = "0.285"
$ws.Range("E45").Value = "  +1.71%  "
$ws.Range("D46").Value = "16.91"
$ws.Range("E46").Value = "  +4.95%  "
$ws.Range("D47").Value = "3.90"
$ws.Range("E47").Value = "  +2.65%  "
$ws.Range("E48").Value = "  -0.51%  "
$ws.Range("D49").Value = "21.66"
$ws.Range("E49").Value = "  -20.83%  "
$ws.Range("D50").Value = "3.769.61"
$ws.Range("E50").Value = "  -7.14%  "
$ws.Range("D51").Value = "2.153.70"
$ws.Range("E51").Value = "  -0.47%  "

$ws.Range("D2:D51").Style = "Normal"
